$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (P4049798I -> P4980809I)
$ws.Range("A2").Value = "SAMSUNG PM P4980809I"
$ws.Range("B2").Value = "P4980809I"

# Update row 3 (P4049798N -> P4980809N)
$ws.Range("A3").Value = "SAMSUNG PM P4980809N"
$ws.Range("B3").Value = "P4980809N"

# Remove rows 4-7 entirely (P4049799I, P4049799N, P4049800N, P4049800I)
$ws.Range("A4:L7").EntireRow.Delete()
